$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record was added at the top of the data block (row 19),
# pushing every existing data row down by one. The former last row (166)
# becomes the new last row (167). Inserting a whole row preserves all the
# existing rows' values/formatting while shifting them down automatically.
$ws.Rows.Item(19).Insert()

# Populate the newly inserted row 19 with this week's record.
$ws.Cells.Item(19, 1).Value = 8
$ws.Cells.Item(19, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(19, 3).Value = "Coquimbo"
$ws.Cells.Item(19, 4).Value = 44490
$ws.Cells.Item(19, 5).Value = 4
$ws.Cells.Item(19, 6).Value = 100112012
$ws.Cells.Item(19, 7).Value = "Espinaca"
$ws.Cells.Item(19, 8).Value = "Sin especificar"
$ws.Cells.Item(19, 9).Value = "Primera"
$ws.Cells.Item(19, 10).Value = 3000
$ws.Cells.Item(19, 11).Value = 400
$ws.Cells.Item(19, 12).Value = 500
$ws.Cells.Item(19, 13).Value = 450
$ws.Cells.Item(19, 14).Value = "`$/atado 300 a 500 gramos"
$ws.Cells.Item(19, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(19, 16).Value = 900
$ws.Cells.Item(19, 17).Value = 0.5
$ws.Cells.Item(19, 18).Value = "Hortaliza"
